$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column D as Text first so numeric-looking price strings
# (e.g. "0.993", "2.637.17") are preserved as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.460.16'
$ws.Range("E2").Value = '  +0.70%  '

$ws.Range("D3").Value = '2.637.17'
$ws.Range("E3").Value = '  -1.06%  '

$ws.Range("D4").Value = '0.993'
$ws.Range("E4").Value = '  -0.66%  '

$ws.Range("D5").Value = '597.13'
$ws.Range("E5").Value = '  +0.37%  '

$ws.Range("D6").Value = '170.13'
$ws.Range("E6").Value = '  +4.04%  '

$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.27%  '

$ws.Range("D8").Value = '0.534'
$ws.Range("E8").Value = '  -1.95%  '

$ws.Range("D9").Value = '2.631.29'
$ws.Range("E9").Value = '  -1.29%  '

$ws.Range("E10").Value = '  +0.51%  '

$ws.Range("E11").Value = '  +1.52%  '

$ws.Range("D12").Value = '0.365'
$ws.Range("E12").Value = '  +2.39%  '

$ws.Range("D13").Value = '5.26'
$ws.Range("E13").Value = '  +1.42%  '

$ws.Range("D14").Value = '27.83'
$ws.Range("E14").Value = '  +0.49%  '

$ws.Range("D15").Value = '3.100.30'
$ws.Range("E15").Value = '  -2.04%  '

$ws.Range("E16").Value = '  +0.58%  '

$ws.Range("D17").Value = '67.181.48'
$ws.Range("E17").Value = '  +0.34%  '

$ws.Range("D18").Value = '2.606.61'
$ws.Range("E18").Value = '  -2.72%  '

$ws.Range("D19").Value = '12.16'
$ws.Range("E19").Value = '  +4.86%  '

$ws.Range("D20").Value = '8.13'
$ws.Range("E20").Value = '  +8.67%  '

$ws.Range("D21").Value = '358.48'
$ws.Range("E21").Value = '  -0.66%  '

$ws.Range("D22").Value = '4.35'
$ws.Range("E22").Value = '  -0.38%  '

$ws.Range("D23").Value = '4.70'
$ws.Range("E23").Value = '  -1.64%  '

$ws.Range("D24").Value = '10.65'
$ws.Range("E24").Value = '  +6.37%  '

$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.19%  '

$ws.Range("B26").Value = 'SuiNetwork'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D26").Value = '1.93'
$ws.Range("E26").Value = '  -4.26%  '

$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value = '70.03'
$ws.Range("E27").Value = '  -1.59%  '

$ws.Range("D28").Value = '2.758.40'

$ws.Range("E29").Value = '  -0.72%  '

$ws.Range("D30").Value = '0.0000101'
$ws.Range("E30").Value = '  -0.34%  '

$ws.Range("D31").Value = '551.79'
$ws.Range("E31").Value = '  -0.19%  '

$ws.Range("D32").Value = '7.96'
$ws.Range("E32").Value = '  +0.21%  '

$ws.Range("D33").Value = '1.36'
$ws.Range("E33").Value = '  -1.43%  '

$ws.Range("D34").Value = '1.91'
$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").Value = '0.137'
$ws.Range("E35").Value = '  +6.82%  '

$ws.Range("D36").Value = '0.991'
$ws.Range("E36").Value = '  -0.80%  '

$ws.Range("D37").Value = '1.51'
$ws.Range("E37").Value = '  -3.67%  '

$ws.Range("D38").Value = '157.62'
$ws.Range("E38").Value = '  +1.95%  '

$ws.Range("D39").Value = '19.05'
$ws.Range("E39").Value = '  -1.81%  '

$ws.Range("D40").Value = '0.368'
$ws.Range("E40").Value = '  -1.09%  '

$ws.Range("D41").Value = '5.22'
$ws.Range("E41").Value = '  -0.84%  '

$ws.Range("D42").Value = '1.81'
$ws.Range("E42").Value = '  -0.37%  '

$ws.Range("D43").Value = '18.13'
$ws.Range("E43").Value = '  +1.31%  '

$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '2.43'
$ws.Range("E45").Value = '  -3.44%  '

$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").Value = '40.12'
$ws.Range("E46").Value = '  -0.29%  '

$ws.Range("E47").Value = '  +0.84%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '0.583'
$ws.Range("E48").Value = '  -0.11%  '

$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '152.28'
$ws.Range("E49").Value = '  -0.08%  '

$ws.Range("D50").Value = '3.80'
$ws.Range("E50").Value = '  -0.46%  '

$ws.Range("D51").Value = '1.71'
$ws.Range("E51").Value = '  -0.42%  '
